$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,15
$data[0,0] = 0
$data[0,1] = 14.98943825577119
$data[0,2] = 8.45486036487449
$data[0,3] = 7.826453265643568
$data[0,4] = 13.18243162950826
$data[0,5] = 39.21268724597837
$data[0,6] = 0
$data[0,7] = 7.344005520526261
$data[0,8] = 0
$data[0,9] = 10.46998192983844
$data[0,10] = 11.17209005334961
$data[0,11] = 10.75886745923933
$data[0,12] = 15.85931230122856
$data[0,13] = 22.14080014276757
$data[0,14] = 30.35078949967897
$data[1,0] = 1
$data[1,1] = 14.83767363624769
$data[1,2] = 8.43378268015935
$data[1,3] = 7.818760874871312
$data[1,4] = 13.20497416956307
$data[1,5] = 39.29629240620882
$data[1,6] = 0
$data[1,7] = 7.344005520526261
$data[1,8] = 0
$data[1,9] = 10.48674579404068
$data[1,10] = 11.05898004020345
$data[1,11] = 10.76701763653773
$data[1,12] = 15.84351969159542
$data[1,13] = 22.20101662380855
$data[1,14] = 30.43016392084462
$data[2,0] = 2
$data[2,1] = 14.74646009688612
$data[2,2] = 8.420710087942245
$data[2,3] = 7.815012378423022
$data[2,4] = 13.22007704451492
$data[2,5] = 39.35443384063958
$data[2,6] = 0
$data[2,7] = 7.344005520526261
$data[2,8] = 0
$data[2,9] = 10.49760597704243
$data[2,10] = 10.99079257564687
$data[2,11] = 10.77309331754195
$data[2,12] = 15.83583246005943
$data[2,13] = 22.23974107025965
$data[2,14] = 30.48369821670014
$data[3,0] = 3
$data[3,1] = 14.70982404545582
$data[3,2] = 8.415351013580977
$data[3,3] = 7.813731502843781
$data[3,4] = 13.22654939544093
$data[3,5] = 39.37983782280882
$data[3,6] = 0
$data[3,7] = 7.344005520526261
$data[3,8] = 0
$data[3,9] = 10.50217458702363
$data[3,10] = 10.96335117137761
$data[3,11] = 10.77583929229319
$data[3,12] = 15.83320834543765
$data[3,13] = 22.25596323037436
$data[3,14] = 30.50672002068894
$data[4,0] = 4
$data[4,1] = 14.70377403356757
$data[4,2] = 8.414459246676227
$data[4,3] = 7.813533756555792
$data[4,4] = 13.22764333519198
$data[4,5] = 39.38415943619461
$data[4,6] = 0
$data[4,7] = 7.344005520526261
$data[4,8] = 0
$data[4,9] = 10.50294184953424
$data[4,10] = 10.95881622507057
$data[4,11] = 10.77631159062841
$data[4,12] = 15.83280341377068
$data[4,13] = 22.25868361835605
$data[4,14] = 30.51061560999844
$data[5,0] = 5
$data[5,1] = 14.74596379689705
$data[5,2] = 8.420637941376333
$data[5,3] = 7.814994103327369
$data[5,4] = 13.22016304538289
$data[5,5] = 39.35476952174257
$data[5,6] = 0
$data[5,7] = 7.344005520526261
$data[5,8] = 0
$data[5,9] = 10.49766701139809
$data[5,10] = 10.99042105556581
$data[5,11] = 10.77312925629362
$data[5,12] = 15.83579500741018
$data[5,13] = 22.23995805815674
$data[5,14] = 30.48400381370398
$data[6,0] = 6
$data[6,1] = 14.93672432002258
$data[6,2] = 8.447620531384858
$data[6,3] = 7.8235996166799
$data[6,4] = 13.18994277144607
$data[6,5] = 39.24010089852857
$data[6,6] = 0
$data[6,7] = 7.344005520526261
$data[6,8] = 0
$data[6,9] = 10.47564464598154
$data[6,10] = 11.13284469143181
$data[6,11] = 10.7614556686006
$data[6,12] = 15.85345191713106
$data[6,13] = 22.16120006803237
$data[6,14] = 30.37716174986323
$data[7,0] = 7
$data[7,1] = 15.32470085877992
$data[7,2] = 8.499452317396218
$data[7,3] = 7.848135364505343
$data[7,4] = 13.14066723275294
$data[7,5] = 39.0692868317745
$data[7,6] = 0
$data[7,7] = 7.344005520526261
$data[7,8] = 0
$data[7,9] = 10.43694055117999
$data[7,10] = 11.42089747315587
$data[7,11] = 10.74703558540016
$data[7,12] = 15.90387347449312
$data[7,13] = 22.02059308063244
$data[7,14] = 30.20573099641368
$data[8,0] = 8
$data[8,1] = 15.61584971034199
$data[8,2] = 8.536817980687843
$data[8,3] = 7.870728217622132
$data[8,4] = 13.11051928720098
$data[8,5] = 38.97678120092689
$data[8,6] = 0
$data[8,7] = 7.344005520526261
$data[8,8] = 0
$data[8,9] = 10.41121161155599
$data[8,10] = 11.63616428011461
$data[8,11] = 10.74156529659041
$data[8,12] = 15.95033816151704
$data[8,13] = 21.92564140667583
$data[8,14] = 30.10301880039714
$data[9,0] = 9
$data[9,1] = 15.74911789993033
$data[9,2] = 8.553649251027252
$data[9,3] = 7.881973576386195
$data[9,4] = 13.09811181105028
$data[9,5] = 38.94186606478488
$data[9,6] = 0
$data[9,7] = 7.344005520526261
$data[9,8] = 0
$data[9,9] = 10.40008929282181
$data[9,10] = 11.73451927196706
$data[9,11] = 10.7401801548923
$data[9,12] = 15.973472227847
$data[9,13] = 21.88424128782879
$data[9,14] = 30.0613414588543
$data[10,0] = 10
$data[10,1] = 15.79965906265225
$data[10,2] = 8.559997775706387
$data[10,3] = 7.88636880944609
$data[10,4] = 13.09360077864886
$data[10,5] = 38.9296749937441
$data[10,6] = 0
$data[10,7] = 7.344005520526261
$data[10,8] = 0
$data[10,9] = 10.3959608350277
$data[10,10] = 11.77179507604466
$data[10,11] = 10.73981344650823
$data[10,12] = 15.98251507315572
$data[10,13] = 21.86882078835313
$data[10,14] = 30.04628526611068
$data[11,0] = 11
$data[11,1] = 15.78877144490825
$data[11,2] = 8.558631642775451
$data[11,3] = 7.885416168297474
$data[11,4] = 13.09456398364745
$data[11,5] = 38.93225472780981
$data[11,6] = 0
$data[11,7] = 7.344005520526261
$data[11,8] = 0
$data[11,9] = 10.3968462720773
$data[11,10] = 11.76376617638731
$data[11,11] = 10.73988541865347
$data[11,12] = 15.98055504857876
$data[11,13] = 21.87213046424692
$data[11,14] = 30.04949558752946
$data[12,0] = 12
$data[12,1] = 15.75327468526612
$data[12,2] = 8.554172057939207
$data[12,3] = 7.882332445551453
$data[12,4] = 13.09773693280491
$data[12,5] = 38.94084244526832
$data[12,6] = 0
$data[12,7] = 7.344005520526261
$data[12,8] = 0
$data[12,9] = 10.39974797421606
$data[12,10] = 11.73758553045733
$data[12,11] = 10.74014682767908
$data[12,12] = 15.9742105546244
$data[12,13] = 21.88296749329284
$data[12,14] = 30.06008822261166
$data[13,0] = 13
$data[13,1] = 15.73154044598998
$data[13,2] = 8.551437128511562
$data[13,3] = 7.880461329338882
$data[13,4] = 13.09970484612838
$data[13,5] = 38.94623687298876
$data[13,6] = 0
$data[13,7] = 7.344005520526261
$data[13,8] = 0
$data[13,9] = 10.40153619120933
$data[13,10] = 11.72155224599977
$data[13,11] = 10.74032747479935
$data[13,12] = 15.97036101242644
$data[13,13] = 21.88963889961058
$data[13,14] = 30.06667108829615
$data[14,0] = 14
$data[14,1] = 15.60715320051837
$data[14,2] = 8.535714566944543
$data[14,3] = 7.870012586366086
$data[14,4] = 13.11135639932385
$data[14,5] = 38.9792071986772
$data[14,6] = 0
$data[14,7] = 7.344005520526261
$data[14,8] = 0
$data[14,9] = 10.41195016029235
$data[14,10] = 11.62974254308151
$data[14,11] = 10.74167794938632
$data[14,12] = 15.9488660852754
$data[14,13] = 21.9283830023025
$data[14,14] = 30.10584410324203
$data[15,0] = 15
$data[15,1] = 15.53102576767076
$data[15,2] = 8.52602580175579
$data[15,3] = 7.863848882916716
$data[15,4] = 13.11883862131269
$data[15,5] = 39.0012689663566
$data[15,6] = 0
$data[15,7] = 7.344005520526261
$data[15,8] = 0
$data[15,9] = 10.41848757964256
$data[15,10] = 11.57350813464886
$data[15,11] = 10.7427884214318
$data[15,12] = 15.93618786601266
$data[15,13] = 21.95260990114662
$data[15,14] = 30.13116840822875
$data[16,0] = 16
$data[16,1] = 15.4873193079104
$data[16,2] = 8.520437473155066
$data[16,3] = 7.860394918275955
$data[16,4] = 13.1232652484797
$data[16,5] = 39.01463283611692
$data[16,6] = 0
$data[16,7] = 7.344005520526261
$data[16,8] = 0
$data[16,9] = 10.42230252295616
$data[16,10] = 11.54120579832063
$data[16,11] = 10.74353103791678
$data[16,12] = 15.9290839567068
$data[16,13] = 21.96671348731497
$data[16,14] = 30.14620930418286
$data[17,0] = 17
$data[17,1] = 15.47253611254334
$data[17,2] = 8.518542715818047
$data[17,3] = 7.859241203193531
$data[17,4] = 13.12478517859282
$data[17,5] = 39.01927345713796
$data[17,6] = 0
$data[17,7] = 7.344005520526261
$data[17,8] = 0
$data[17,9] = 10.42360361874854
$data[17,10] = 11.53027697551964
$data[17,11] = 10.74380034546465
$data[17,12] = 15.92671116951392
$data[17,13] = 21.97151776514173
$data[17,14] = 30.15138346686851
$data[18,0] = 18
$data[18,1] = 15.53912169315865
$data[18,2] = 8.527058811017792
$data[18,3] = 7.864495593578569
$data[18,4] = 13.11802939453035
$data[18,5] = 38.99885064154724
$data[18,6] = 0
$data[18,7] = 7.344005520526261
$data[18,8] = 0
$data[18,9] = 10.41778599194068
$data[18,10] = 11.57949023925845
$data[18,11] = 10.74265946226731
$data[18,12] = 15.93751803246266
$data[18,13] = 21.9500134325771
$data[18,14] = 30.12842342647613
$data[19,0] = 19
$data[19,1] = 15.76369923798984
$data[19,2] = 8.555482636939532
$data[19,3] = 7.883234514221021
$data[19,4] = 13.09679987875783
$data[19,5] = 38.93829205702324
$data[19,6] = 0
$data[19,7] = 7.344005520526261
$data[19,8] = 0
$data[19,9] = 10.39889341540593
$data[19,10] = 11.74527482557385
$data[19,11] = 10.74006576935139
$data[19,12] = 15.97606645850618
$data[19,13] = 21.87977743390709
$data[19,14] = 30.05695719989169
$data[20,0] = 20
$data[20,1] = 15.91089353748761
$data[20,2] = 8.573912333320106
$data[20,3] = 7.896278062267472
$data[20,4] = 13.08401727732524
$data[20,5] = 38.9047200912793
$data[20,6] = 0
$data[20,7] = 7.344005520526261
$data[20,8] = 0
$data[20,9] = 10.38703153230939
$data[20,10] = 11.85379078331003
$data[20,11] = 10.73929008997701
$data[20,12] = 16.00290423633168
$data[20,13] = 21.83537061650434
$data[20,14] = 30.01448234336405
$data[21,0] = 21
$data[21,1] = 15.83230881152295
$data[21,2] = 8.564089886189869
$data[21,3] = 7.889244376313161
$data[21,4] = 13.09073984131769
$data[21,5] = 38.92208853217491
$data[21,6] = 0
$data[21,7] = 7.344005520526261
$data[21,8] = 0
$data[21,9] = 10.39331813472529
$data[21,10] = 11.79586863825757
$data[21,11] = 10.73962024239205
$data[21,12] = 15.98843157849138
$data[21,13] = 21.85893481425647
$data[21,14] = 30.03676460190052
$data[22,0] = 22
$data[22,1] = 15.5354613343973
$data[22,2] = 8.526591843916544
$data[22,3] = 7.86420293633681
$data[22,4] = 13.11839485656913
$data[22,5] = 38.99994184722067
$data[22,6] = 0
$data[22,7] = 7.344005520526261
$data[22,8] = 0
$data[22,9] = 10.41810300374094
$data[22,10] = 11.57678564151831
$data[22,11] = 10.74271744010117
$data[22,12] = 15.93691608770252
$data[22,13] = 21.95118675012468
$data[22,14] = 30.12966293254692
$data[23,0] = 23
$data[23,1] = 15.21849856602434
$data[23,2] = 8.485551925346831
$data[23,3] = 7.840687938428345
$data[23,4] = 13.15293186810355
$data[23,5] = 39.10970493415715
$data[23,6] = 0
$data[23,7] = 7.344005520526261
$data[23,8] = 0
$data[23,9] = 10.43694055117999
$data[23,10] = 11.34220918585589
$data[23,11] = 10.75003368283881
$data[23,12] = 15.88856272663742
$data[23,13] = 22.0571582745776
$data[23,14] = 30.24802790656703

$range = $ws.Range("A2:O25")
$range.Value = $data
Write-Host "done"